# ajustes no modelo e remoção da parte de setores
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "[-ENDEREÇO 3-]" paragraph: split off a new empty right-aligned
#    paragraph right after it (mirrors the paragraph-mark formatting
#    that the diff adds to the original paragraph's pPr/rPr).
# ---------------------------------------------------------------------
$found = $d.Content.Find.Execute("[-ENDEREÇO 3-]", $false, $false, $false, $false, $false, $true, 1, $false, "[-ENDEREÇO 3-]^p", 2)
if ($found) {
    $i = 0
    foreach ($p in $d.Paragraphs) {
        $i = $i + 1
        if ($p.Range.Text -eq "[-ENDEREÇO 3-]`r") {
            $r = $p.Range
            $r.Font.Name = "Calibri"
            $r.Font.Color = 0
            break
        }
    }
}

# ---------------------------------------------------------------------
# 2) Remove the (now redundant) blank right-aligned paragraph that used
#    to sit right after "[-MEMBRO 4-]".
# ---------------------------------------------------------------------
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    $t = $p.Range.Text
    if ($t -eq "[-MEMBRO 4-]`r") {
        $nextP = $d.Paragraphs($i + 1)
        if ($nextP.Range.Text -eq "`r") {
            $nextP.Range.Delete()
        }
        break
    }
}

# ---------------------------------------------------------------------
# 3) "... pontos de fuga nos setores [-SETORES-]." -> "... pontos de
#    fuga na instalação." (drop the highlighted [-SETORES-] placeholder)
# ---------------------------------------------------------------------
$d.Content.Find.Execute(" pontos de fuga nos setores [-SETORES-]", $false, $false, $false, $false, $false, $true, 1, $false, " pontos de fuga na instalação", 2) | Out-Null

# ---------------------------------------------------------------------
# 4) "... de acordo com o setor e local em que foram identificados." ->
#    "... de acordo com o local em que foram identificados."
# ---------------------------------------------------------------------
$oldTail = "seguir é apresentada a relação de pontos de fuga encontrados de acordo com o setor e local em que foram identificados."
$newTail = "seguir é apresentada a relação de pontos de fuga encontrados de acordo com o local em que foram identificados."
$d.Content.Find.Execute($oldTail, $false, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null
